$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# This text lives in one shared string, referenced by every cell that showed
# "Ready for handoff" -- the Overview rollup columns as well as each language
# sheet's Status column -- so all of them flip to the new text together.
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime ---
# zh-cn: both rows handed back at 2016-03-18 18:14:53
$wsZh.Range("H2").Value = "2016-03-18 18:14:53"
$wsZh.Range("H3").Value = "2016-03-18 18:14:53"
# de-de: both rows handed back at 2016-03-18 18:14:58
$wsDe.Range("H2").Value = "2016-03-18 18:14:58"
$wsDe.Range("H3").Value = "2016-03-18 18:14:58"

# --- New columns: F = Latest Target File, G = Latest Handback File ---
# (the handback file matches the already-handed-off target file 1:1, since
# the translation came back "in sync" -- same filename/link as column D)
# zh-cn row 2 (4d2394fc...)
$wsZh.Range("F2").Value = "4d2394fc-73e6-419f-83a2-02da558e7b5a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/21b16df405d549cda272a2f5fa70fccb724e612f/e2e/4d2394fc-73e6-419f-83a2-02da558e7b5a.md", "", "", "4d2394fc-73e6-419f-83a2-02da558e7b5a.md") | Out-Null

$wsZh.Range("G2").Value = "4d2394fc-73e6-419f-83a2-02da558e7b5a.efa79e593f7824f447a50b06f94124913ed4274e.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4011018654639dfbe8363d7daa5151a6a2531692/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/4d2394fc-73e6-419f-83a2-02da558e7b5a.efa79e593f7824f447a50b06f94124913ed4274e.zh-cn.xlf", "", "", "4d2394fc-73e6-419f-83a2-02da558e7b5a.efa79e593f7824f447a50b06f94124913ed4274e.zh-cn.xlf") | Out-Null

# zh-cn row 3 (e0e28513...)
$wsZh.Range("F3").Value = "e0e28513-fdec-42e5-8632-bd52df9df165.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/21b16df405d549cda272a2f5fa70fccb724e612f/e2e/e0e28513-fdec-42e5-8632-bd52df9df165.md", "", "", "e0e28513-fdec-42e5-8632-bd52df9df165.md") | Out-Null

$wsZh.Range("G3").Value = "e0e28513-fdec-42e5-8632-bd52df9df165.11f7a4e51c55b514c6602583b671f381b8ad4d1a.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4011018654639dfbe8363d7daa5151a6a2531692/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/e0e28513-fdec-42e5-8632-bd52df9df165.11f7a4e51c55b514c6602583b671f381b8ad4d1a.zh-cn.xlf", "", "", "e0e28513-fdec-42e5-8632-bd52df9df165.11f7a4e51c55b514c6602583b671f381b8ad4d1a.zh-cn.xlf") | Out-Null

# de-de row 2 (4d2394fc...)
$wsDe.Range("F2").Value = "4d2394fc-73e6-419f-83a2-02da558e7b5a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/21b16df405d549cda272a2f5fa70fccb724e612f/e2e/4d2394fc-73e6-419f-83a2-02da558e7b5a.md", "", "", "4d2394fc-73e6-419f-83a2-02da558e7b5a.md") | Out-Null

$wsDe.Range("G2").Value = "4d2394fc-73e6-419f-83a2-02da558e7b5a.efa79e593f7824f447a50b06f94124913ed4274e.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aace33db0a95e7e9c39bf873a44686656cf74e6d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/4d2394fc-73e6-419f-83a2-02da558e7b5a.efa79e593f7824f447a50b06f94124913ed4274e.de-de.xlf", "", "", "4d2394fc-73e6-419f-83a2-02da558e7b5a.efa79e593f7824f447a50b06f94124913ed4274e.de-de.xlf") | Out-Null

# de-de row 3 (e0e28513...)
$wsDe.Range("F3").Value = "e0e28513-fdec-42e5-8632-bd52df9df165.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/21b16df405d549cda272a2f5fa70fccb724e612f/e2e/e0e28513-fdec-42e5-8632-bd52df9df165.md", "", "", "e0e28513-fdec-42e5-8632-bd52df9df165.md") | Out-Null

$wsDe.Range("G3").Value = "e0e28513-fdec-42e5-8632-bd52df9df165.11f7a4e51c55b514c6602583b671f381b8ad4d1a.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aace33db0a95e7e9c39bf873a44686656cf74e6d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/e0e28513-fdec-42e5-8632-bd52df9df165.11f7a4e51c55b514c6602583b671f381b8ad4d1a.de-de.xlf", "", "", "e0e28513-fdec-42e5-8632-bd52df9df165.11f7a4e51c55b514c6602583b671f381b8ad4d1a.de-de.xlf") | Out-Null

"Handback report generated"
